$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (column 9), shifting existing I..M to J..N
$ws.Columns.Item(9).Insert()

# New column I (rows 2-66) gets the style used by column G/H (hyperlink-like
# style, s="2" in the OOXML), matching the diff.
$ws.Range("G2:G66").Copy()
$ws.Range("I2:I66").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New column I formula (rows 2 through 66), per the commit message
# ("started text decode") - a tweaked variant of the hex formatting formula.
$ws.Range("I2").Formula = '=G2 &"=0x" & RIGHT("0" &  DEC2HEX( MID(C2,5,10)),2) & ","'
$ws.Range("I3:I66").FormulaR1C1 = $ws.Range("I2").FormulaR1C1

# Selection as captured in the diff
$ws.Range("I2:I66").Select()
$ws.Range("I66").Activate()
